# Reorder the "Recorded By" (column G) comma-separated list so that any
# entry containing "System" alongside other recorder name(s) lists
# "System" first, reversing the original order of the list.
# Single-value cells (just one name, with or without "System") are left
# untouched, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text

    if ($null -eq $current -or $current -eq "") {
        continue
    }

    $parts = $current -split ", "

    if ($parts.Count -gt 1 -and ($parts -contains "System")) {
        $reversed = $parts[($parts.Count - 1)..0]
        $newValue = $reversed -join ", "
        $cell.Value = $newValue
    }
}
